# Fruta / hortaliza, semanal
# Insert a new weekly price record at row 564 (pushing existing rows 564-641
# down to 565-642) and update the sheet dimension accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 564; all rows below shift down by one.
$ws.Range("A564:R564").EntireRow.Insert()

# Populate the newly inserted row with the new data point.
$ws.Range("A564").Value = 3
$ws.Range("B564").Value = "Femacal de La Calera"
$ws.Range("C564").Value = "Coquimbo"
$ws.Range("D564").Value = 45131
$ws.Range("E564").Value = 5
$ws.Range("F564").Value = 100112031
$ws.Range("G564").Value = "Poroto verde"
$ws.Range("H564").Value = "Magnum"
$ws.Range("I564").Value = "Primera"
$ws.Range("J564").Value = 80
$ws.Range("K564").Value = 26000
$ws.Range("L564").Value = 27000
$ws.Range("M564").Value = 26500
$ws.Range("N564").Value = "`$/malla 25 kilos"
$ws.Range("O564").Value = "Región de Arica y Parinacota"
$ws.Range("P564").Value = 1060
$ws.Range("Q564").Value = 25
$ws.Range("R564").Value = "Hortaliza"

# Match the date number format used by the rest of column D.
$ws.Range("D564").NumberFormat = $ws.Range("D565").NumberFormat
